$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in columns A and B (rows 1-32)
$ws.Cells.Item(1, 1).Value = -0.22018249879080543
$ws.Cells.Item(1, 2).Value = 0.21958288703395112
$ws.Cells.Item(2, 1).Value = -0.14852849596766227
$ws.Cells.Item(2, 2).Value = 0.14684735605327592
$ws.Cells.Item(3, 1).Value = -0.11827330717625451
$ws.Cells.Item(3, 2).Value = 0.11780281601642706
$ws.Cells.Item(4, 1).Value = -0.10980281607207942
$ws.Cells.Item(4, 2).Value = 0.10938180461776881
$ws.Cells.Item(5, 1).Value = -0.10638180464924929
$ws.Cells.Item(5, 2).Value = 0.10495436754385334
$ws.Cells.Item(6, 1).Value = -0.020741866330732606
$ws.Cells.Item(6, 2).Value = 0.020605645912500137
$ws.Cells.Item(7, 1).Value = -0.010605645990318546
$ws.Cells.Item(7, 2).Value = 0.010586530089320867
$ws.Cells.Item(8, 1).Value = -0.00058653016757315157
$ws.Cells.Item(8, 2).Value = 0.0005778249279182468
$ws.Cells.Item(9, 1).Value = -0.024174493046965839
$ws.Cells.Item(9, 2).Value = 0.024030334446965895
$ws.Cells.Item(10, 1).Value = -0.022030334482170844
$ws.Cells.Item(10, 2).Value = 0.022022066080966241
$ws.Cells.Item(11, 1).Value = -0.019022066121957337
$ws.Cells.Item(11, 2).Value = 0.019008009173873575
$ws.Cells.Item(12, 1).Value = -0.015508009218133889
$ws.Cells.Item(12, 2).Value = 0.015412164683787477
$ws.Cells.Item(13, 1).Value = -0.011912164729607433
$ws.Cells.Item(13, 2).Value = 0.011874723960577427
$ws.Cells.Item(14, 1).Value = -0.0038747240315100129
$ws.Cells.Item(14, 2).Value = 0.0038685762850025895
$ws.Cells.Item(15, 1).Value = -0.0028685763172413559
$ws.Cells.Item(15, 2).Value = 0.002866826651882981
$ws.Cells.Item(16, 1).Value = -0.0060350607524064515
$ws.Cells.Item(16, 2).Value = 0.0060034916562980989
$ws.Cells.Item(17, 1).Value = -0.0040034916945028698
$ws.Cells.Item(17, 2).Value = 0.0039999999505528905
$ws.Cells.Item(18, 1).Value = -0.065054391100936471
$ws.Cells.Item(18, 2).Value = 0.064927334806267112
$ws.Cells.Item(19, 1).Value = -0.06092733483071866
$ws.Cells.Item(19, 2).Value = 0.059983598221077461
$ws.Cells.Item(20, 1).Value = -0.0080170041261240499
$ws.Cells.Item(20, 2).Value = 0.0080057617157205385
$ws.Cells.Item(21, 1).Value = -0.0040057617493989284
$ws.Cells.Item(21, 2).Value = 0.0039999999661128882
$ws.Cells.Item(22, 1).Value = -0.045718708458034385
$ws.Cells.Item(22, 2).Value = 0.045503431329488642
$ws.Cells.Item(23, 1).Value = -0.040503431367008957
$ws.Cells.Item(23, 2).Value = 0.040099897617377422
$ws.Cells.Item(24, 1).Value = -0.020099897741101103
$ws.Cells.Item(24, 2).Value = 0.01999999987465273
$ws.Cells.Item(25, 1).Value = -0.082212501247013847
$ws.Cells.Item(25, 2).Value = 0.082124628625175333
$ws.Cells.Item(26, 1).Value = -0.079624628663125918
$ws.Cells.Item(26, 2).Value = 0.079513015195164627
$ws.Cells.Item(27, 1).Value = -0.077013015235116278
$ws.Cells.Item(27, 2).Value = 0.076361370058938949
$ws.Cells.Item(28, 1).Value = -0.074361370102826285
$ws.Cells.Item(28, 2).Value = 0.073923420247743898
$ws.Cells.Item(29, 1).Value = -0.06692342032513654
$ws.Cells.Item(29, 2).Value = 0.066803210094043664
$ws.Cells.Item(30, 1).Value = -0.0068032104675759264
$ws.Cells.Item(30, 2).Value = 0.0067518954809950849
$ws.Cells.Item(31, 1).Value = -0.014023625424606223
$ws.Cells.Item(31, 2).Value = 0.014001137283020881
$ws.Cells.Item(32, 1).Value = -0.0040011373809569761
$ws.Cells.Item(32, 2).Value = 0.003999999935198062

# Column B width changed from 16.42578125 to 15.7109375 (characters)
# ColumnWidth rounds to the nearest achievable pixel width in this engine;
# 14.8 is the COM value that lands on the closest achievable stored width.
$ws.Columns.Item(2).ColumnWidth = 14.8
